$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "Public" column (C) to TRUE for rows 15 through 32
$ws.Range("C15:C32").Value = $true

# Update the active selection to match the range that was checked (C15:C38)
$ws.Range("C15:C38").Select()
